# Gendata.xlsx edit:
#  - Update x_ohm_per_km (column F) values on the "Lines" sheet for rows 2 and 3
#    from 0.0083 to 0.083.
#  - Make "Lines" the active/selected sheet (previously "Gen slack" was active),
#    and set its selection to E9 (previously J10).

$wb = $excel.ActiveWorkbook

$wsLines = $wb.Worksheets.Item("Lines")

# Update the x_ohm_per_km values for the two line entries.
$wsLines.Range("F2").Value = 0.083
$wsLines.Range("F3").Value = 0.083

# Switch the active sheet to "Lines" (this also clears tabSelected on the
# previously active "Gen slack" sheet) and set the new selection.
$wsLines.Activate() | Out-Null
$wsLines.Range("E9").Select() | Out-Null
